$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 5556055.5
$ws.Range("I33").Value = 7143425
$ws.Range("J33").Value = 262.5
$ws.Range("K33").Value = 7143425
$ws.Range("L33").Value = 262.5
$ws.Range("M33").Value = -7143196
$ws.Range("N33").Value = -720.5
$ws.Range("H43").Value = 13055.556
$ws.Range("I43").Value = 17625
$ws.Range("J43").Value = 9400
$ws.Range("K43").Value = 17625
$ws.Range("L43").Value = 9400
$ws.Range("M43").Value = -17556
$ws.Range("N43").Value = -9538
$ws.Range("H51").Value = 18886.736
$ws.Range("I51").Value = 8499.333000000001
$ws.Range("J51").Value = 20834.375
$ws.Range("K51").Value = 8499.333000000001
$ws.Range("L51").Value = 20834.375
$ws.Range("M51").Value = -8015.333000000001
$ws.Range("N51").Value = -21802.375
$ws.Range("H62").Value = 74512690
$ws.Range("I62").Value = 97438530
$ws.Range("K62").Value = 97438530
$ws.Range("M62").Value = -97437906
$ws.Range("H65").Value = 74512690
$ws.Range("I65").Value = 97438530
$ws.Range("K65").Value = 487192650
$ws.Range("M65").Value = -487189530
$ws.Range("H70").Value = 5672.231
$ws.Range("J70").Value = 5875.4
$ws.Range("L70").Value = 17626.2
$ws.Range("N70").Value = -18166.2
$ws.Range("H73").Value = 5672.231
$ws.Range("J73").Value = 5875.4
$ws.Range("L73").Value = 17626.2
$ws.Range("N73").Value = -19498.2
$ws.Range("H100").Value = 21921.176
$ws.Range("I100").Value = 56591.11
$ws.Range("J100").Value = 9440
$ws.Range("K100").Value = 56591.11
$ws.Range("L100").Value = 9440
$ws.Range("M100").Value = -56050.11
$ws.Range("N100").Value = -10522
$ws.Range("H105").Value = 15000
$ws.Range("J105").Value = 15000
$ws.Range("L105").Value = 15000
$ws.Range("N105").Value = -21988
$ws.Range("H107").Value = 15625794
$ws.Range("I107").Value = 17857912
$ws.Range("K107").Value = 17857912
$ws.Range("M107").Value = -17855992
$ws.Range("H108").Value = 19592
$ws.Range("J108").Value = 19592
$ws.Range("L108").Value = 19592
$ws.Range("N108").Value = -27272
$ws.Range("H112").Value = 2849.3076
$ws.Range("J112").Value = 3021.3914
$ws.Range("L112").Value = 9064.174199999999
$ws.Range("N112").Value = -11280.1742
$ws.Range("H116").Value = 8033
$ws.Range("J116").Value = 8333
$ws.Range("L116").Value = 8333
$ws.Range("N116").Value = -15217

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 503110.8
$ws.Range("I74").Value = 909774.25
$ws.Range("K74").Value = 909774.25
$ws.Range("M74").Value = -908900.25
$ws.Range("H77").Value = 503110.8
$ws.Range("I77").Value = 909774.25
$ws.Range("K77").Value = 4548871.25
$ws.Range("M77").Value = -4544503.25
$ws.Range("H122").Value = 1500.1786
$ws.Range("I122").Value = 1042.3334
$ws.Range("K122").Value = 3127.0002
$ws.Range("M122").Value = -677.0001999999999

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 90932690
$ws.Range("I105").Value = 100025870
$ws.Range("K105").Value = 100025870
$ws.Range("M105").Value = -100024123

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 47623108
$ws.Range("I31").Value = 111112620
$ws.Range("J31").Value = 5970.1665
$ws.Range("K31").Value = 111112620
$ws.Range("L31").Value = 5970.1665
$ws.Range("M31").Value = -111112325
$ws.Range("N31").Value = -6560.1665
$ws.Range("H34").Value = 47623108
$ws.Range("I34").Value = 111112620
$ws.Range("J34").Value = 5970.1665
$ws.Range("K34").Value = 111112620
$ws.Range("L34").Value = 5970.1665
$ws.Range("M34").Value = -111112418
$ws.Range("N34").Value = -6374.1665
$ws.Range("H107").Value = 1325.6786
$ws.Range("I107").Value = 1055.76
$ws.Range("J107").Value = 3575
$ws.Range("K107").Value = 1055.76
$ws.Range("L107").Value = 3575
$ws.Range("M107").Value = 864.24
$ws.Range("N107").Value = -7415
$ws.Range("H132").Value = 112046.63
$ws.Range("I132").Value = 3314.125
$ws.Range("J132").Value = 402000
$ws.Range("K132").Value = 9942.375
$ws.Range("L132").Value = 1206000
$ws.Range("M132").Value = -7412.375
$ws.Range("N132").Value = -1211060
$ws.Range("H134").Value = 8299.576999999999
$ws.Range("I134").Value = 8464.817999999999
$ws.Range("K134").Value = 25394.454
$ws.Range("M134").Value = -22859.454

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 3246.0454
$ws.Range("I2").Value = 117.42857
$ws.Range("K2").Value = 704.57142
$ws.Range("M2").Value = -591.57142
$ws.Range("H80").Value = 4094.5
$ws.Range("J80").Value = 4094.5
$ws.Range("L80").Value = 12283.5
$ws.Range("N80").Value = -14155.5
$ws.Range("H83").Value = 4094.5
$ws.Range("J83").Value = 4094.5
$ws.Range("L83").Value = 36850.5
$ws.Range("N83").Value = -46210.5
$ws.Range("H103").Value = 3012.3845
$ws.Range("J103").Value = 2718.5
$ws.Range("L103").Value = 8155.5
$ws.Range("N103").Value = -9913.5
$ws.Range("H132").Value = 1749.7646
$ws.Range("I132").Value = 1308.8
$ws.Range("J132").Value = 2379.7144
$ws.Range("K132").Value = 11779.2
$ws.Range("L132").Value = 21417.4296
$ws.Range("M132").Value = -9249.199999999999
$ws.Range("N132").Value = -26477.4296

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 832.625
$ws.Range("I2").Value = 1733
$ws.Range("J2").Value = 292.4
$ws.Range("K2").Value = 1733
$ws.Range("L2").Value = 292.4
$ws.Range("M2").Value = -1620
$ws.Range("N2").Value = -518.4
$ws.Range("H126").Value = 3011.5
$ws.Range("I126").Value = 2410.6667
$ws.Range("K126").Value = 7232.000100000001
$ws.Range("M126").Value = -4762.000100000001

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3244.1353
$ws.Range("J22").Value = 4603.8237
$ws.Range("L22").Value = 4603.8237
$ws.Range("N22").Value = -5193.8237
$ws.Range("H27").Value = 3244.1353
$ws.Range("J27").Value = 4603.8237
$ws.Range("L27").Value = 4603.8237
$ws.Range("N27").Value = -4817.8237
$ws.Range("H46").Value = 4612.4
$ws.Range("I46").Value = 1345.7142
$ws.Range("K46").Value = 1345.7142
$ws.Range("M46").Value = -1157.7142
$ws.Range("H68").Value = 2299
$ws.Range("J68").Value = 2400
$ws.Range("L68").Value = 2400
$ws.Range("N68").Value = -3898
$ws.Range("H71").Value = 2299
$ws.Range("J71").Value = 2400
$ws.Range("L71").Value = 12000
$ws.Range("N71").Value = -19488
$ws.Range("H93").Value = 3012.037
$ws.Range("I93").Value = 2238.125
$ws.Range("K93").Value = 2238.125
$ws.Range("M93").Value = -990.125
$ws.Range("H122").Value = 3571.182
$ws.Range("I122").Value = 3528.3
$ws.Range("K122").Value = 10584.9
$ws.Range("M122").Value = -8134.900000000001

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()
$ws.Range("H28").Value = 43500
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 43500
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 43500
$ws.Range("M28").ClearContents()
$ws.Range("N28").Value = -44196
$ws.Range("H31").Value = 58179.332
$ws.Range("J31").Value = 58179.332
$ws.Range("L31").Value = 58179.332
$ws.Range("N31").Value = -58875.332
$ws.Range("H33").Value = 29340
$ws.Range("J33").Value = 29340
$ws.Range("L33").Value = 29340
$ws.Range("N33").Value = -29840
$ws.Range("H36").Value = 29340
$ws.Range("J36").Value = 29340
$ws.Range("L36").Value = 29340
$ws.Range("N36").Value = -29840
$ws.Range("H126").Value = 13499.5
$ws.Range("I126").Value = 12000
$ws.Range("K126").Value = 36000
$ws.Range("M126").Value = -33530
$ws.Range("H132").Value = 4485.587
$ws.Range("I132").Value = 2469.6428
$ws.Range("J132").Value = 7621.5
$ws.Range("K132").Value = 7408.928400000001
$ws.Range("L132").Value = 22864.5
$ws.Range("M132").Value = -4878.928400000001
$ws.Range("N132").Value = -27924.5
